$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1838.3334  # H111: 1916.6666 -> 1838.3334
$ws.Cells.Item(111, 9).Value = 1838.3334  # I111: 1916.6666 -> 1838.3334
$ws.Cells.Item(111, 11).Value = 5515.0002  # K111: 5749.9998 -> 5515.0002
$ws.Cells.Item(111, 13).Value = -2448.0002  # M111: -2682.9998 -> -2448.0002
$ws.Cells.Item(127, 8).Value = 797.1667  # H127: 754.7143 -> 797.1667
$ws.Cells.Item(127, 10).Value = 1369  # J127: 1151.75 -> 1369
$ws.Cells.Item(127, 12).Value = 4107  # L127: 3455.25 -> 4107
$ws.Cells.Item(127, 14).Value = -14027  # N127: -13375.25 -> -14027
$ws.Cells.Item(130, 8).Value = 333373730  # H130: 250041330 -> 333373730
$ws.Cells.Item(130, 10).Value = 333373730  # J130: 250041330 -> 333373730
$ws.Cells.Item(130, 12).Value = 333373730  # L130: 250041330 -> 333373730
$ws.Cells.Item(130, 14).Value = -333383770  # N130: -250051370 -> -333383770
$ws.Cells.Item(138, 8).Value = 2179.625  # H138: 2079.9214 -> 2179.625
$ws.Cells.Item(138, 9).Value = 1368.9584  # I138: 1325.0392 -> 1368.9584
$ws.Cells.Item(138, 10).Value = 2990.2917  # J138: 3093.0527 -> 2990.2917
$ws.Cells.Item(138, 11).Value = 4106.8752  # K138: 3975.1176 -> 4106.8752
$ws.Cells.Item(138, 12).Value = 8970.875100000001  # L138: 9279.158100000001 -> 8970.875100000001
$ws.Cells.Item(138, 13).Value = 1033.1248  # M138: 1164.8824 -> 1033.1248
$ws.Cells.Item(138, 14).Value = -19250.8751  # N138: -19559.1581 -> -19250.8751

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(37, 8).Value = 19785.715  # H37: 19000 -> 19785.715
$ws.Cells.Item(37, 9).Value = 18900  # I37: 19000 -> 18900
$ws.Cells.Item(37, 10).Value = 22000  # J37: 0 -> 22000
$ws.Cells.Item(37, 11).Value = 18900  # K37: 19000 -> 18900
$ws.Cells.Item(37, 12).Value = 22000  # L37: 0 -> 22000
$ws.Cells.Item(37, 13).Value = -18627  # M37: -18727 -> -18627
$ws.Cells.Item(37, 14).Value = -22546  # N37: None -> -22546
$ws.Cells.Item(61, 8).Value = 1417.1364  # H61: 2274.7576 -> 1417.1364
$ws.Cells.Item(61, 9).Value = 1090.5385  # I61: 2823.3076 -> 1090.5385
$ws.Cells.Item(61, 10).Value = 1888.8889  # J61: 1918.2 -> 1888.8889
$ws.Cells.Item(61, 11).Value = 1090.5385  # K61: 2823.3076 -> 1090.5385
$ws.Cells.Item(61, 12).Value = 1888.8889  # L61: 1918.2 -> 1888.8889
$ws.Cells.Item(61, 13).Value = -878.5385000000001  # M61: -2611.3076 -> -878.5385000000001
$ws.Cells.Item(61, 14).Value = -2312.8889  # N61: -2342.2 -> -2312.8889
$ws.Cells.Item(97, 8).Value = 878.8333  # H97: 1128.1765 -> 878.8333
$ws.Cells.Item(97, 9).Value = 736.2857  # I97: 833.4545000000001 -> 736.2857
$ws.Cells.Item(97, 10).Value = 1377.75  # J97: 1668.5 -> 1377.75
$ws.Cells.Item(97, 11).Value = 736.2857  # K97: 833.4545000000001 -> 736.2857
$ws.Cells.Item(97, 12).Value = 1377.75  # L97: 1668.5 -> 1377.75
$ws.Cells.Item(97, 13).Value = -240.2857  # M97: -337.4545000000001 -> -240.2857
$ws.Cells.Item(97, 14).Value = -2369.75  # N97: -2660.5 -> -2369.75
$ws.Cells.Item(102, 8).Value = 2033.48  # H102: 2162.45 -> 2033.48
$ws.Cells.Item(102, 9).Value = 1472.6111  # I102: 1593.2667 -> 1472.6111
$ws.Cells.Item(102, 10).Value = 3475.7144  # J102: 3870 -> 3475.7144
$ws.Cells.Item(102, 11).Value = 1472.6111  # K102: 1593.2667 -> 1472.6111
$ws.Cells.Item(102, 12).Value = 3475.7144  # L102: 3870 -> 3475.7144
$ws.Cells.Item(102, 13).Value = 149.3888999999999  # M102: 28.7333000000001 -> 149.3888999999999
$ws.Cells.Item(102, 14).Value = -6719.7144  # N102: -7114 -> -6719.7144
$ws.Cells.Item(105, 8).Value = 38800  # H105: 37134.5 -> 38800
$ws.Cells.Item(105, 10).Value = 38800  # J105: 37134.5 -> 38800
$ws.Cells.Item(105, 12).Value = 38800  # L105: 37134.5 -> 38800
$ws.Cells.Item(105, 14).Value = -45788  # N105: -44122.5 -> -45788
$ws.Cells.Item(132, 8).Value = 3132740  # H132: 2666510 -> 3132740
$ws.Cells.Item(132, 9).Value = 7744.9546  # I132: 11362.071 -> 7744.9546
$ws.Cells.Item(132, 10).Value = 6952178.5  # J132: 3792936.2 -> 6952178.5
$ws.Cells.Item(132, 11).Value = 23234.8638  # K132: 34086.213 -> 23234.8638
$ws.Cells.Item(132, 12).Value = 20856535.5  # L132: 11378808.6 -> 20856535.5
$ws.Cells.Item(132, 13).Value = -20704.8638  # M132: -31556.213 -> -20704.8638
$ws.Cells.Item(132, 14).Value = -20861595.5  # N132: -11383868.6 -> -20861595.5
$ws.Cells.Item(136, 8).Value = 1417.1364  # H136: 2274.7576 -> 1417.1364
$ws.Cells.Item(136, 9).Value = 1090.5385  # I136: 2823.3076 -> 1090.5385
$ws.Cells.Item(136, 10).Value = 1888.8889  # J136: 1918.2 -> 1888.8889
$ws.Cells.Item(136, 11).Value = 3271.6155  # K136: 8469.9228 -> 3271.6155
$ws.Cells.Item(136, 12).Value = 5666.6667  # L136: 5754.6 -> 5666.6667
$ws.Cells.Item(136, 13).Value = -721.6155000000003  # M136: -5919.9228 -> -721.6155000000003
$ws.Cells.Item(136, 14).Value = -10766.6667  # N136: -10854.6 -> -10766.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 13215.286  # H82: 16206.4 -> 13215.286
$ws.Cells.Item(82, 9).Value = 3001.4  # I82: 3201.5557 -> 3001.4
$ws.Cells.Item(82, 10).Value = 38750  # J82: 35713.668 -> 38750
$ws.Cells.Item(82, 11).Value = 3001.4  # K82: 3201.5557 -> 3001.4
$ws.Cells.Item(82, 12).Value = 38750  # L82: 35713.668 -> 38750
$ws.Cells.Item(82, 13).Value = -2618.4  # M82: -2818.5557 -> -2618.4
$ws.Cells.Item(82, 14).Value = -39516  # N82: -36479.668 -> -39516
$ws.Cells.Item(85, 8).Value = 13215.286  # H85: 16206.4 -> 13215.286
$ws.Cells.Item(85, 9).Value = 3001.4  # I85: 3201.5557 -> 3001.4
$ws.Cells.Item(85, 10).Value = 38750  # J85: 35713.668 -> 38750
$ws.Cells.Item(85, 11).Value = 3001.4  # K85: 3201.5557 -> 3001.4
$ws.Cells.Item(85, 12).Value = 38750  # L85: 35713.668 -> 38750
$ws.Cells.Item(85, 13).Value = -1675.4  # M85: -1875.5557 -> -1675.4
$ws.Cells.Item(85, 14).Value = -41402  # N85: -38365.668 -> -41402
$ws.Cells.Item(105, 8).Value = 2264.3333  # H105: 2635.182 -> 2264.3333
$ws.Cells.Item(105, 9).Value = 2111.1428  # I105: 2733.3333 -> 2111.1428
$ws.Cells.Item(105, 10).Value = 2398.375  # J105: 2598.375 -> 2398.375
$ws.Cells.Item(105, 11).Value = 2111.1428  # K105: 2733.3333 -> 2111.1428
$ws.Cells.Item(105, 12).Value = 2398.375  # L105: 2598.375 -> 2398.375
$ws.Cells.Item(105, 13).Value = -364.1428000000001  # M105: -986.3332999999998 -> -364.1428000000001
$ws.Cells.Item(105, 14).Value = -5892.375  # N105: -6092.375 -> -5892.375
$ws.Cells.Item(109, 8).Value = 29745  # H109: 30195 -> 29745
$ws.Cells.Item(109, 10).Value = 29745  # J109: 30195 -> 29745
$ws.Cells.Item(109, 12).Value = 29745  # L109: 30195 -> 29745
$ws.Cells.Item(109, 14).Value = -32519  # N109: -32969 -> -32519
$ws.Cells.Item(134, 8).Value = 2563.487  # H134: 2023.2727 -> 2563.487
$ws.Cells.Item(134, 9).Value = 1414  # I134: 1121.1428 -> 1414
$ws.Cells.Item(134, 10).Value = 4051.0588  # J134: 3041.8064 -> 4051.0588
$ws.Cells.Item(134, 11).Value = 4242  # K134: 3363.4284 -> 4242
$ws.Cells.Item(134, 12).Value = 12153.1764  # L134: 9125.4192 -> 12153.1764
$ws.Cells.Item(134, 13).Value = -1707  # M134: -828.4284000000002 -> -1707
$ws.Cells.Item(134, 14).Value = -17223.1764  # N134: -14195.4192 -> -17223.1764

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 3758.5925  # H7: 3493.138 -> 3758.5925
$ws.Cells.Item(7, 9).Value = 12552.625  # I7: 14301 -> 12552.625
$ws.Cells.Item(7, 10).Value = 55.842106  # J7: 54.272728 -> 55.842106
$ws.Cells.Item(7, 11).Value = 12552.625  # K7: 14301 -> 12552.625
$ws.Cells.Item(7, 12).Value = 55.842106  # L7: 54.272728 -> 55.842106
$ws.Cells.Item(7, 13).Value = -12439.625  # M7: -14188 -> -12439.625
$ws.Cells.Item(7, 14).Value = -281.842106  # N7: -280.272728 -> -281.842106
$ws.Cells.Item(22, 8).Value = 763.5333000000001  # H22: 791.82355 -> 763.5333000000001
$ws.Cells.Item(22, 9).Value = 395.25  # I22: 411.57144 -> 395.25
$ws.Cells.Item(22, 10).Value = 897.4545000000001  # J22: 1058 -> 897.4545000000001
$ws.Cells.Item(22, 11).Value = 395.25  # K22: 411.57144 -> 395.25
$ws.Cells.Item(22, 12).Value = 897.4545000000001  # L22: 1058 -> 897.4545000000001
$ws.Cells.Item(22, 13).Value = -45.25  # M22: -61.57144 -> -45.25
$ws.Cells.Item(22, 14).Value = -1597.4545  # N22: -1758 -> -1597.4545
$ws.Cells.Item(107, 8).Value = 1333.619  # H107: 1228.75 -> 1333.619
$ws.Cells.Item(107, 9).Value = 388  # I107: 382.41666 -> 388
$ws.Cells.Item(107, 10).Value = 2193.2727  # J107: 2498.25 -> 2193.2727
$ws.Cells.Item(107, 11).Value = 388  # K107: 382.41666 -> 388
$ws.Cells.Item(107, 12).Value = 2193.2727  # L107: 2498.25 -> 2193.2727
$ws.Cells.Item(107, 13).Value = 1532  # M107: 1537.58334 -> 1532
$ws.Cells.Item(107, 14).Value = -6033.2727  # N107: -6338.25 -> -6033.2727
$ws.Cells.Item(114, 8).Value = 37950  # H114: 30680 -> 37950
$ws.Cells.Item(114, 10).Value = 37950  # J114: 30680 -> 37950
$ws.Cells.Item(114, 12).Value = 37950  # L114: 30680 -> 37950
$ws.Cells.Item(114, 14).Value = -46628  # N114: -39358 -> -46628
$ws.Cells.Item(123, 8).Value = 19800  # H123: 21018 -> 19800
$ws.Cells.Item(123, 9).Value = 0  # I123: 15000 -> 0
$ws.Cells.Item(123, 10).Value = 19800  # J123: 21686.666 -> 19800
$ws.Cells.Item(123, 11).Value = 0  # K123: 15000 -> 0
$ws.Cells.Item(123, 12).Value = 19800  # L123: 21686.666 -> 19800
$ws.Cells.Item(123, 13).ClearContents()  # M123: -10100 -> (removed)
$ws.Cells.Item(123, 14).Value = -29600  # N123: -31486.666 -> -29600
$ws.Cells.Item(134, 8).Value = 1964.2  # H134: 2060.6191 -> 1964.2
$ws.Cells.Item(134, 9).Value = 1425.7273  # I134: 1636.75 -> 1425.7273
$ws.Cells.Item(134, 10).Value = 3445  # J134: 2321.4614 -> 3445
$ws.Cells.Item(134, 11).Value = 4277.1819  # K134: 4910.25 -> 4277.1819
$ws.Cells.Item(134, 12).Value = 10335  # L134: 6964.3842 -> 10335
$ws.Cells.Item(134, 13).Value = -1742.1819  # M134: -2375.25 -> -1742.1819
$ws.Cells.Item(134, 14).Value = -15405  # N134: -12034.3842 -> -15405

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 19255.363  # H55: 1201.4286 -> 19255.363
$ws.Cells.Item(55, 9).Value = 971.2857  # I55: 500 -> 971.2857
$ws.Cells.Item(55, 10).Value = 51252.5  # J55: 1318.3334 -> 51252.5
$ws.Cells.Item(55, 11).Value = 2913.8571  # K55: 1500 -> 2913.8571
$ws.Cells.Item(55, 12).Value = 153757.5  # L55: 3955.0002 -> 153757.5
$ws.Cells.Item(55, 13).Value = -2736.8571  # M55: -1323 -> -2736.8571
$ws.Cells.Item(55, 14).Value = -154111.5  # N55: -4309.0002 -> -154111.5
$ws.Cells.Item(122, 8).Value = 1809.127  # H122: 1721.3088 -> 1809.127
$ws.Cells.Item(122, 9).Value = 519.7692  # I122: 466.52942 -> 519.7692
$ws.Cells.Item(122, 10).Value = 2144.36  # J122: 2139.5686 -> 2144.36
$ws.Cells.Item(122, 11).Value = 4677.922799999999  # K122: 4198.76478 -> 4677.922799999999
$ws.Cells.Item(122, 12).Value = 19299.24  # L122: 19256.1174 -> 19299.24
$ws.Cells.Item(122, 13).Value = -2227.922799999999  # M122: -1748.76478 -> -2227.922799999999
$ws.Cells.Item(122, 14).Value = -24199.24  # N122: -24156.1174 -> -24199.24

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(42, 8).Value = 0  # H42: 40000 -> 0
$ws.Cells.Item(42, 10).Value = 0  # J42: 40000 -> 0
$ws.Cells.Item(42, 12).Value = 0  # L42: 40000 -> 0
$ws.Cells.Item(42, 14).ClearContents()  # N42: -40970 -> (removed)
$ws.Cells.Item(113, 8).Value = 1822.2142  # H113: 1763.1875 -> 1822.2142
$ws.Cells.Item(113, 9).Value = 1690.1111  # I113: 1611.1 -> 1690.1111
$ws.Cells.Item(113, 10).Value = 2060  # J113: 2016.6666 -> 2060
$ws.Cells.Item(113, 11).Value = 1690.1111  # K113: 1611.1 -> 1690.1111
$ws.Cells.Item(113, 12).Value = 2060  # L113: 2016.6666 -> 2060
$ws.Cells.Item(113, 13).Value = 479.8888999999999  # M113: 558.9000000000001 -> 479.8888999999999
$ws.Cells.Item(113, 14).Value = -6400  # N113: -6356.6666 -> -6400
$ws.Cells.Item(115, 8).Value = 0  # H115: 40000 -> 0
$ws.Cells.Item(115, 10).Value = 0  # J115: 40000 -> 0
$ws.Cells.Item(115, 12).Value = 0  # L115: 40000 -> 0
$ws.Cells.Item(115, 14).ClearContents()  # N115: -42350 -> (removed)
$ws.Cells.Item(123, 8).Value = 18589.732  # H123: 18857.066 -> 18589.732
$ws.Cells.Item(123, 10).Value = 18589.732  # J123: 18857.066 -> 18589.732
$ws.Cells.Item(123, 12).Value = 18589.732  # L123: 18857.066 -> 18589.732
$ws.Cells.Item(123, 14).Value = -23489.732  # N123: -23757.066 -> -23489.732

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(29, 8).Value = 0  # H29: 5000 -> 0
$ws.Cells.Item(29, 10).Value = 0  # J29: 5000 -> 0
$ws.Cells.Item(29, 12).Value = 0  # L29: 5000 -> 0
$ws.Cells.Item(29, 14).Value = 0  # N29: -5590 -> 0
$ws.Cells.Item(106, 8).Value = 17342.5  # H106: 20185 -> 17342.5
$ws.Cells.Item(106, 10).Value = 17342.5  # J106: 20185 -> 17342.5
$ws.Cells.Item(106, 12).Value = 17342.5  # L106: 20185 -> 17342.5
$ws.Cells.Item(106, 14).Value = -19866.5  # N106: -22709 -> -19866.5
$ws.Cells.Item(108, 8).Value = 21750  # H108: 24666.666 -> 21750
$ws.Cells.Item(108, 10).Value = 21750  # J108: 24666.666 -> 21750
$ws.Cells.Item(108, 12).Value = 21750  # L108: 24666.666 -> 21750
$ws.Cells.Item(108, 14).Value = -29430  # N108: -32346.666 -> -29430
$ws.Cells.Item(136, 8).Value = 1563.2667  # H136: 1718.2667 -> 1563.2667
$ws.Cells.Item(136, 9).Value = 1323.3438  # I136: 1397.4839 -> 1323.3438
$ws.Cells.Item(136, 10).Value = 2153.8462  # J136: 2428.5715 -> 2153.8462
$ws.Cells.Item(136, 11).Value = 3970.0314  # K136: 4192.4517 -> 3970.0314
$ws.Cells.Item(136, 12).Value = 6461.5386  # L136: 7285.7145 -> 6461.5386
$ws.Cells.Item(136, 13).Value = -1420.0314  # M136: -1642.4517 -> -1420.0314
$ws.Cells.Item(136, 14).Value = -11561.5386  # N136: -12385.7145 -> -11561.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 9092.5  # H104: 9423.333000000001 -> 9092.5
$ws.Cells.Item(104, 10).Value = 9092.5  # J104: 9423.333000000001 -> 9092.5
$ws.Cells.Item(104, 12).Value = 9092.5  # L104: 9423.333000000001 -> 9092.5
$ws.Cells.Item(104, 14).Value = -16080.5  # N104: -16411.333 -> -16080.5
$ws.Cells.Item(119, 8).Value = 56750  # H119: 40017.816 -> 56750
$ws.Cells.Item(119, 10).Value = 56750  # J119: 40017.816 -> 56750
$ws.Cells.Item(119, 12).Value = 56750  # L119: 40017.816 -> 56750
$ws.Cells.Item(119, 14).Value = -66426  # N119: -49693.816 -> -66426
